$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Shared-string text fixup: every label ending in ": NoSettings" becomes
#    ": test". All of these labels live in column A of "BAU Emissions"
#    (rows 4-280), so iterate over just that used range instead of the
#    whole sheet to minimise unrelated churn.
# ---------------------------------------------------------------------------
$emisWs = $wb.Worksheets.Item("BAU Emissions")
$labelRange = $emisWs.Range("A4:A280")
$labelRange.Replace(": NoSettings", ": test")

# ---------------------------------------------------------------------------
# 2) "About" sheet: bump the report date in C1 (was 3/18/2024 -> 4/5/2024).
# ---------------------------------------------------------------------------
$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("C1").Value = 45387

# ---------------------------------------------------------------------------
# 3) Row 94 ("...natural gas if,iron and steel 241,CO2] : test") on
#    "BAU Emissions" gets new figures for 2032-2050 (columns M:AE).
# ---------------------------------------------------------------------------
$newValues = @(1001080, 2002150, 3003230, 4004300, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380)
$col = 13
foreach ($val in $newValues) {
    $emisWs.Cells.Item(94, $col).Value = $val
    $col = $col + 1
}

# ---------------------------------------------------------------------------
# 4) Update the saved selection on "BAU Emissions" (was activeCell B283).
# ---------------------------------------------------------------------------
$emisWs.Activate()
$emisWs.Range("A30:AE280").Select()

# ---------------------------------------------------------------------------
# 5) Make "About" the active/selected sheet again (was "Current and Planned
#    Capacity"), keeping its existing E29 selection.
# ---------------------------------------------------------------------------
$aboutWs.Activate()
$aboutWs.Range("E29").Select()
